$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 1034, shifting existing rows 1034+ down by one
$ws.Rows.Item(1034).Insert()

# Populate the newly inserted row 1034 with the new data
$ws.Cells.Item(1034, 1).Value = 3
$ws.Cells.Item(1034, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1034, 3).Value = "Coquimbo"
$ws.Cells.Item(1034, 4).Value = 45223
$ws.Cells.Item(1034, 5).Value = 5
$ws.Cells.Item(1034, 6).Value = 100112006
$ws.Cells.Item(1034, 7).Value = "Repollo"
$ws.Cells.Item(1034, 8).Value = "Crespo record"
$ws.Cells.Item(1034, 9).Value = "Primera"
$ws.Cells.Item(1034, 10).Value = 2300
$ws.Cells.Item(1034, 11).Value = 800
$ws.Cells.Item(1034, 12).Value = 900
$ws.Cells.Item(1034, 13).Value = 852
$ws.Cells.Item(1034, 14).Value = "`$/unidad"
$ws.Cells.Item(1034, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1034, 16).Value = 852
$ws.Cells.Item(1034, 17).Value = 1
$ws.Cells.Item(1034, 18).Value = "Hortaliza"

# Apply the date style (s="2") to the new D1034 cell, matching other date cells in column D
$ws.Cells.Item(1034, 4).NumberFormat = $ws.Cells.Item(1035, 4).NumberFormat
